$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$singles = @(
    "Tử Vi tọa thủ tại cung đối Thiên Di",
    "Thiên Cơ tọa thủ tại cung đối Thiên Di",
    "Thái Dương tọa thủ tại cung đối Thiên Di",
    "Vũ Khúc tọa thủ tại cung đối Thiên Di",
    "Thiên Đồng tọa thủ tại cung đối Thiên Di",
    "Liêm Trinh tọa thủ tại cung đối Thiên Di",
    "Thiên Phủ tọa thủ tại cung đối Thiên Di",
    "Thái Âm tọa thủ tại cung đối Thiên Di",
    "Tham Lang tọa thủ tại cung đối Thiên Di",
    "Cự Môn tọa thủ tại cung đối Thiên Di",
    "Thiên Tướng tọa thủ tại cung đối Thiên Di",
    "Thiên Lương tọa thủ tại cung đối Thiên Di",
    "Thất Sát tọa thủ tại cung đối Thiên Di",
    "Phá Quân tọa thủ tại cung đối Thiên Di"
)

$pairs = @(
    "Tử Vi đồng cung Thiên Cơ tại cung đối Thiên Di",
    "Tử Vi đồng cung Thái Dương tại cung đối Thiên Di",
    "Tử Vi đồng cung Vũ Khúc tại cung đối Thiên Di",
    "Tử Vi đồng cung Thiên Đồng tại cung đối Thiên Di",
    "Tử Vi đồng cung Liêm Trinh tại cung đối Thiên Di",
    "Tử Vi đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Tử Vi đồng cung Thái Âm tại cung đối Thiên Di",
    "Tử Vi đồng cung Tham Lang tại cung đối Thiên Di",
    "Tử Vi đồng cung Cự Môn tại cung đối Thiên Di",
    "Tử Vi đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Tử Vi đồng cung Thiên Lương tại cung đối Thiên Di",
    "Tử Vi đồng cung Thất Sát tại cung đối Thiên Di",
    "Tử Vi đồng cung Phá Quân tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thái Dương tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Vũ Khúc tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thiên Đồng tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Liêm Trinh tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thái Âm tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Tham Lang tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Cự Môn tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Thất Sát tại cung đối Thiên Di",
    "Thiên Cơ đồng cung Phá Quân tại cung đối Thiên Di",
    "Thái Dương đồng cung Vũ Khúc tại cung đối Thiên Di",
    "Thái Dương đồng cung Thiên Đồng tại cung đối Thiên Di",
    "Thái Dương đồng cung Liêm Trinh tại cung đối Thiên Di",
    "Thái Dương đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Thái Dương đồng cung Thái Âm tại cung đối Thiên Di",
    "Thái Dương đồng cung Tham Lang tại cung đối Thiên Di",
    "Thái Dương đồng cung Cự Môn tại cung đối Thiên Di",
    "Thái Dương đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Thái Dương đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thái Dương đồng cung Thất Sát tại cung đối Thiên Di",
    "Thái Dương đồng cung Phá Quân tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thiên Đồng tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Liêm Trinh tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thái Âm tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Tham Lang tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Cự Môn tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thiên Lương tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Thất Sát tại cung đối Thiên Di",
    "Vũ Khúc đồng cung Phá Quân tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Liêm Trinh tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Thái Âm tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Tham Lang tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Cự Môn tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Thất Sát tại cung đối Thiên Di",
    "Thiên Đồng đồng cung Phá Quân tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Thiên Phủ tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Thái Âm tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Tham Lang tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Cự Môn tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Thiên Lương tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Thất Sát tại cung đối Thiên Di",
    "Liêm Trinh đồng cung Phá Quân tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Thái Âm tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Tham Lang tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Cự Môn tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Thất Sát tại cung đối Thiên Di",
    "Thiên Phủ đồng cung Phá Quân tại cung đối Thiên Di",
    "Thái Âm đồng cung Tham Lang tại cung đối Thiên Di",
    "Thái Âm đồng cung Cự Môn tại cung đối Thiên Di",
    "Thái Âm đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Thái Âm đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thái Âm đồng cung Thất Sát tại cung đối Thiên Di",
    "Thái Âm đồng cung Phá Quân tại cung đối Thiên Di",
    "Tham Lang đồng cung Cự Môn tại cung đối Thiên Di",
    "Tham Lang đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Tham Lang đồng cung Thiên Lương tại cung đối Thiên Di",
    "Tham Lang đồng cung Thất Sát tại cung đối Thiên Di",
    "Tham Lang đồng cung Phá Quân tại cung đối Thiên Di",
    "Cự Môn đồng cung Thiên Tướng tại cung đối Thiên Di",
    "Cự Môn đồng cung Thiên Lương tại cung đối Thiên Di",
    "Cự Môn đồng cung Thất Sát tại cung đối Thiên Di",
    "Cự Môn đồng cung Phá Quân tại cung đối Thiên Di",
    "Thiên Tướng đồng cung Thiên Lương tại cung đối Thiên Di",
    "Thiên Tướng đồng cung Thất Sát tại cung đối Thiên Di",
    "Thiên Tướng đồng cung Phá Quân tại cung đối Thiên Di",
    "Thiên Lương đồng cung Thất Sát tại cung đối Thiên Di",
    "Thiên Lương đồng cung Phá Quân tại cung đối Thiên Di",
    "Thất Sát đồng cung Phá Quân tại cung đối Thiên Di"
)

$startRow = 4307
for ($i = 0; $i -lt $singles.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $singles[$i]
    $ws.Cells.Item($row, 2).Value = $singles[$i]
}

$pairStartRow = 4322
for ($i = 0; $i -lt $pairs.Length; $i++) {
    $row = $pairStartRow + $i
    $ws.Cells.Item($row, 1).Value = $pairs[$i]
    $ws.Cells.Item($row, 2).Value = $pairs[$i]
}

$ws.Range("B4322:B4412").Select()
